$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ur = $ws.UsedRange
$lastRow = $ur.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()

    if ($val -ne $null) {
        $text = $val.ToString()

        if ($text.Contains(",")) {
            $rawParts = $text.Split(",")
            $parts = @()
            foreach ($p in $rawParts) {
                $parts += $p.Trim()
            }

            if ($parts[0].Equals("System")) {
                $others = @()
                $systems = @()
                foreach ($p in $parts) {
                    if ($p.Equals("System")) {
                        $systems += $p
                    } else {
                        $others += $p
                    }
                }
                $newParts = $others + $systems
                $newText = $newParts -join ", "
                $cell.Value = $newText
            }
        }
    }
}
